# AR TCs updated for Document Number field change
# Insert a new "DivisionNumber" column (column C) into the CashReceipt and
# CashReceipt_ForeignCurr sheets, right after the "Deposit Amount Bank"
# column, shifting the former "CompanyNumber" column (and everything to its
# right) one column over. The new column mirrors the value already present
# in the (now-shifted) CompanyNumber column.

$wb = $excel.ActiveWorkbook

# ---- CashReceipt ----------------------------------------------------
$ws9 = $wb.Worksheets.Item("CashReceipt")
$null = $ws9.Columns("C:C").Insert()
$ws9.Range("C1").Value = "DivisionNumber"
$ws9.Range("C2").Value = 100
$ws9.Columns("C:C").ColumnWidth = 17.5

# ---- CashReceipt_ForeignCurr -----------------------------------------
$ws10 = $wb.Worksheets.Item("CashReceipt_ForeignCurr")
$null = $ws10.Columns("C:C").Insert()
$ws10.Range("C1").Value = "DivisionNumber"
$ws10.Range("C2").Value = 100
$ws10.Columns("C:C").ColumnWidth = 18.5

# ---- restore on-screen selections ------------------------------------
# CashReceipt's selection moved to C5; CashReceipt_ForeignCurr (the tab
# that stays active) moved to E8 -- select it last so it remains the
# active sheet/tab.
$null = $ws9.Range("C5").Select()
$null = $ws10.Range("E8").Select()
